# Rerun mode choice for nonGerman cities:
#  - Rename the summary sheets to reflect new run identifiers
#  - Update the "CarAvailable" variable label to "CarOwnershipHH" on every sheet

$wb = $excel.ActiveWorkbook

# Map old sheet names -> new sheet names (order matches tab/rId order)
$renameMap = @{
    "summ6"  = "summ32266441"
    "summ0"  = "summ25900342"
    "summ9"  = "summ18786055"
    "summ15" = "summ16911695"
    "summ13" = "summ12144344"
}

foreach ($ws in $wb.Worksheets) {
    if ($renameMap.ContainsKey($ws.Name)) {
        $ws.Name = $renameMap[$ws.Name]
    }
}

# Update the variable label in column A, row 31 on every worksheet
foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("A31")
    if ($cell.Text -eq "CarAvailable") {
        $cell.Value = "CarOwnershipHH"
    }
}
